$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: fill in previously-empty keep/reason/block columns ---
$ws.Range("I17").Value = "keep"
$ws.Range("J17").Value = "NA"
$ws.Range("K17").Value = "NA"
$ws.Range("H17").Font.Color = 4342338

# --- Row 18: new participant SAN-080618-01 ---
$ws.Range("A18").Value = "SAN-080618-01"
$ws.Range("B18").Value = "IF"
$ws.Range("C18").Value = (Get-Date -Year 2018 -Month 8 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D18").Value = "M"
$ws.Range("E18").Value = "C"
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = "nogaze_gaze_ol_gr"
$ws.Range("H18").Font.Color = 4342338
$ws.Range("I18").Value = "keep"
$ws.Range("J18").Value = "NA"
$ws.Range("K18").Value = "NA"
$ws.Range("L18").Value = "Wore glasses but eye tracker seemed mostly fine"

# --- Row 19: new participant SAN-080818-01 ---
$ws.Range("A19").Value = "SAN-080818-01"
$ws.Range("B19").Value = "WM"
$ws.Range("C19").Value = (Get-Date -Year 2018 -Month 8 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D19").Value = "M"
$ws.Range("E19").Value = "C"
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = "nogaze_gaze_gr_ol"
$ws.Range("H19").Font.Color = 4342338
$ws.Range("I19").Value = "keep"
$ws.Range("J19").Value = "NA"
$ws.Range("K19").Value = "NA"
$ws.Range("L19").Value = "Liked when speaker smiled; asked for copy of results"

# --- Row 20: new participant SAN-081018-01 ---
$ws.Range("A20").Value = "SAN-081018-01"
$ws.Range("B20").Value = "AL"
$ws.Range("C20").Value = (Get-Date -Year 2018 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D20").Value = "M"
$ws.Range("E20").Value = "A"
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = "nogaze_gaze_ol_gr"
$ws.Range("H20").Font.Color = 4342338
$ws.Range("I20").Value = "keep"
$ws.Range("J20").Value = "NA"
$ws.Range("K20").Value = "NA"
$ws.Range("L20").Value = "Eye tracker was fine during calibration but seemed to have some trouble later in experiment"

# --- Update selection ---
$ws.Range("C22").Select()
